$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.514509333333333
$ws.Range("H2").Value = 4.543528
$ws.Range("I2").Value = 0.01996786707219448
$ws.Range("J2").Value = 0.02165084619119693
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 184.662213217064
$ws.Range("R2").Value = 1661.959918953576
$ws.Range("S2").Value = 0.004557130822921913
$ws.Range("T2").Value = 0.005238899488419402
$ws.Range("G3").Value = 1.514509333333333
$ws.Range("H3").Value = 4.543528
$ws.Range("I3").Value = 0.01996786707219448
$ws.Range("J3").Value = 0.02165084619119693
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 224.0176333187466
$ws.Range("R3").Value = 2016.15869986872
$ws.Range("S3").Value = 0.005528351707097073
$ws.Range("T3").Value = 0.006355419683024021
$ws.Range("G4").Value = 1.514509333333333
$ws.Range("H4").Value = 4.543528
$ws.Range("I4").Value = 0.01996786707219448
$ws.Range("J4").Value = 0.02165084619119693
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 126.4690554350471
$ws.Range("R4").Value = 1138.221498915424
$ws.Range("S4").Value = 0.003121028501870117
$ws.Range("T4").Value = 0.003587949360493902
$ws.Range("G5").Value = 1.514509333333333
$ws.Range("H5").Value = 4.543528
$ws.Range("I5").Value = 0.01996786707219448
$ws.Range("J5").Value = 0.02165084619119693
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 137.9240364285667
$ws.Range("R5").Value = 827.5442185714002
$ws.Range("S5").Value = 0.003403716800965673
$ws.Range("T5").Value = 0.002608619458192388
$ws.Range("G6").Value = 1.514509333333333
$ws.Range("H6").Value = 4.543528
$ws.Range("I6").Value = 0.01996786707219448
$ws.Range("J6").Value = 0.02165084619119693
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 136.05690010088
$ws.Range("R6").Value = 1224.51210090792
$ws.Range("S6").Value = 0.003357639239339698
$ws.Range("T6").Value = 0.003859958201067221
$ws.Range("I7").Value = 0.6527104067845205
$ws.Range("J7").Value = 0.7077236929508544
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 6036.245527419434
$ws.Range("R7").Value = 54326.2097467749
$ws.Range("S7").Value = 0.1489636675988119
$ws.Range("T7").Value = 0.1712493479562032
$ws.Range("I8").Value = 0.6527104067845205
$ws.Range("J8").Value = 0.7077236929508544
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("Q8").Value = 7322.697013242646
$ws.Range("S8").Value = 0.180710973212156
$ws.Range("T8").Value = 0.2077462030168184
$ws.Range("I9").Value = 0.6527104067845205
$ws.Range("J9").Value = 0.7077236929508544
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 4134.025347835594
$ws.Range("R9").Value = 37206.22813052035
$ws.Range("S9").Value = 0.1020202996983316
$ws.Range("T9").Value = 0.117283026681971
$ws.Range("I10").Value = 0.6527104067845205
$ws.Range("J10").Value = 0.7077236929508544
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 4508.466207089942
$ws.Range("R10").Value = 27050.79724253965
$ws.Range("S10").Value = 0.1112608256908559
$ws.Range("T10").Value = 0.08527065317225414
$ws.Range("I11").Value = 0.6527104067845205
$ws.Range("J11").Value = 0.7077236929508544
$ws.Range("M11").Value = 89.83562999999999
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 4447.43318300378
$ws.Range("R11").Value = 40026.89864703402
$ws.Range("S11").Value = 0.109754640584365
$ws.Range("T11").Value = 0.1261744621236075
$ws.Range("G12").Value = 3.794695333333333
$ws.Range("H12").Value = 11.384086
$ws.Range("I12").Value = 0.05003070653167101
$ws.Range("J12").Value = 0.05424751316892035
$ws.Range("M12").Value = 121.928739
$ws.Range("N12").Value = 365.786217
$ws.Range("O12").Value = 0.2282232151508951
$ws.Range("P12").Value = 0.2419720431319445
$ws.Range("Q12").Value = 462.682416882518
$ws.Range("R12").Value = 4164.141751942661
$ws.Range("S12").Value = 0.01141816870092884
$ws.Range("T12").Value = 0.01312638159631072
$ws.Range("G13").Value = 3.794695333333333
$ws.Range("H13").Value = 11.384086
$ws.Range("I13").Value = 0.05003070653167101
$ws.Range("J13").Value = 0.05424751316892035
$ws.Range("O13").Value = 0.2768624053389947
$ws.Range("P13").Value = 0.2935413991166814
$ws.Range("Q13").Value = 561.2898177841266
$ws.Range("R13").Value = 5051.608360057139
$ws.Range("S13").Value = 0.01385162175116779
$ws.Range("T13").Value = 0.01592389091420548
$ws.Range("G14").Value = 3.794695333333333
$ws.Range("H14").Value = 11.384086
$ws.Range("I14").Value = 0.05003070653167101
$ws.Range("J14").Value = 0.05424751316892035
$ws.Range("M14").Value = 83.50496933333334
$ws.Range("N14").Value = 250.514908
$ws.Range("O14").Value = 0.1563025480180701
$ws.Range("P14").Value = 0.1657186665504434
$ws.Range("Q14").Value = 316.8759174393431
$ws.Range("R14").Value = 2851.883256954088
$ws.Range("S14").Value = 0.007819926910044479
$ws.Range("T14").Value = 0.008989825546031097
$ws.Range("G15").Value = 3.794695333333333
$ws.Range("H15").Value = 11.384086
$ws.Range("I15").Value = 0.05003070653167101
$ws.Range("J15").Value = 0.05424751316892035
$ws.Range("M15").Value = 91.06846250000001
$ws.Range("N15").Value = 182.136925
$ws.Range("O15").Value = 0.1704597085236707
$ws.Range("P15").Value = 0.1204857969594293
$ws.Range("Q15").Value = 345.5770696625917
$ws.Range("R15").Value = 2073.46241797555
$ws.Range("S15").Value = 0.00852821965262195
$ws.Range("T15").Value = 0.006536054857224506
$ws.Range("G16").Value = 3.794695333333333
$ws.Range("H16").Value = 11.384086
$ws.Range("I16").Value = 0.05003070653167101
$ws.Range("J16").Value = 0.05424751316892035
$ws.Range("M16").Value = 89.83562999999999
$ws.Range("N16").Value = 269.50689
$ws.Range("O16").Value = 0.1681521229683693
$ws.Range("P16").Value = 0.1782820942415013
$ws.Range("Q16").Value = 340.89884592806
$ws.Range("R16").Value = 3068.08961335254
$ws.Range("S16").Value = 0.008412769516907941
$ws.Range("T16").Value = 0.00967136025514854
$ws.Range("G17").Value = 17.6874565
$ws.Range("H17").Value = 35.374913
$ws.Range("I17").Value = 0.2331981536620147
$ws.Range("J17").Value = 0.1685687422615142
$ws.Range("M17").Value = 121.928739
$ws.Range("N17").Value = 365.786217
$ws.Range("O17").Value = 0.2282232151508951
$ws.Range("P17").Value = 0.2419720431319445
$ws.Range("Q17").Value = 2156.609267162353
$ws.Range("R17").Value = 12939.65560297412
$ws.Range("S17").Value = 0.05322123239599748
$ws.Range("T17").Value = 0.04078892297320075
$ws.Range("G18").Value = 17.6874565
$ws.Range("H18").Value = 35.374913
$ws.Range("I18").Value = 0.2331981536620147
$ws.Range("J18").Value = 0.1685687422615142
$ws.Range("O18").Value = 0.2768624053389947
$ws.Range("P18").Value = 0.2935413991166814
$ws.Range("Q18").Value = 2616.228277601645
$ws.Range("R18").Value = 15697.36966560987
$ws.Range("S18").Value = 0.06456380174347789
$ws.Range("T18").Value = 0.04948190445078413
$ws.Range("G19").Value = 17.6874565
$ws.Range("H19").Value = 35.374913
$ws.Range("I19").Value = 0.2331981536620147
$ws.Range("J19").Value = 0.1685687422615142
$ws.Range("M19").Value = 83.50496933333334
$ws.Range("N19").Value = 250.514908
$ws.Range("O19").Value = 0.1563025480180701
$ws.Range("P19").Value = 0.1657186665504434
$ws.Range("Q19").Value = 1476.990512617167
$ws.Range("R19").Value = 8861.943075703004
$ws.Range("S19").Value = 0.03644946561048233
$ws.Range("T19").Value = 0.0279349871896635
$ws.Range("G20").Value = 17.6874565
$ws.Range("H20").Value = 35.374913
$ws.Range("I20").Value = 0.2331981536620147
$ws.Range("J20").Value = 0.1685687422615142
$ws.Range("M20").Value = 91.06846250000001
$ws.Range("N20").Value = 182.136925
$ws.Range("O20").Value = 0.1704597085236707
$ws.Range("P20").Value = 0.1204857969594293
$ws.Range("Q20").Value = 1610.769468990631
$ws.Range("R20").Value = 6443.077875962525
$ws.Range("S20").Value = 0.03975088930148521
$ws.Range("T20").Value = 0.02031013925382717
$ws.Range("G21").Value = 17.6874565
$ws.Range("H21").Value = 35.374913
$ws.Range("I21").Value = 0.2331981536620147
$ws.Range("J21").Value = 0.1685687422615142
$ws.Range("M21").Value = 89.83562999999999
$ws.Range("N21").Value = 269.50689
$ws.Range("O21").Value = 0.1681521229683693
$ws.Range("P21").Value = 0.1782820942415013
$ws.Range("Q21").Value = 1588.963797775095
$ws.Range("R21").Value = 9533.78278665057
$ws.Range("S21").Value = 0.03921276461057178
$ws.Range("T21").Value = 0.03005278839403861
$ws.Range("G22").Value = 3.344326
$ws.Range("H22").Value = 10.032978
$ws.Range("I22").Value = 0.04409286594959943
$ws.Range("J22").Value = 0.04780920542751418
$ws.Range("M22").Value = 121.928739
$ws.Range("N22").Value = 365.786217
$ws.Range("O22").Value = 0.2282232151508951
$ws.Range("P22").Value = 0.2419720431319445
$ws.Range("Q22").Value = 407.769451984914
$ws.Range("R22").Value = 3669.925067864226
$ws.Range("S22").Value = 0.01006301563223501
$ws.Range("T22").Value = 0.01156849111781046
$ws.Range("G23").Value = 3.344326
$ws.Range("H23").Value = 10.032978
$ws.Range("I23").Value = 0.04409286594959943
$ws.Range("J23").Value = 0.04780920542751418
$ws.Range("O23").Value = 0.2768624053389947
$ws.Range("P23").Value = 0.2935413991166814
$ws.Range("Q23").Value = 494.6737395915799
$ws.Range("R23").Value = 4452.06365632422
$ws.Range("S23").Value = 0.01220765692509595
$ws.Range("T23").Value = 0.01403398105184935
$ws.Range("G24").Value = 3.344326
$ws.Range("H24").Value = 10.032978
$ws.Range("I24").Value = 0.04409286594959943
$ws.Range("J24").Value = 0.04780920542751418
$ws.Range("M24").Value = 83.50496933333334
$ws.Range("N24").Value = 250.514908
$ws.Range("O24").Value = 0.1563025480180701
$ws.Range("P24").Value = 0.1657186665504434
$ws.Range("Q24").Value = 279.2678400706694
$ws.Range("R24").Value = 2513.410560636024
$ws.Range("S24").Value = 0.006891827297341591
$ws.Range("T24").Value = 0.007922877772283869
$ws.Range("G25").Value = 3.344326
$ws.Range("H25").Value = 10.032978
$ws.Range("I25").Value = 0.04409286594959943
$ws.Range("J25").Value = 0.04780920542751418
$ws.Range("M25").Value = 91.06846250000001
$ws.Range("N25").Value = 182.136925
$ws.Range("O25").Value = 0.1704597085236707
$ws.Range("P25").Value = 0.1204857969594293
$ws.Range("Q25").Value = 304.562626918775
$ws.Range("R25").Value = 1827.37576151265
$ws.Range("S25").Value = 0.007516057077742005
$ws.Range("T25").Value = 0.00576033021793112
$ws.Range("G26").Value = 3.344326
$ws.Range("H26").Value = 10.032978
$ws.Range("I26").Value = 0.04409286594959943
$ws.Range("J26").Value = 0.04780920542751418
$ws.Range("M26").Value = 89.83562999999999
$ws.Range("N26").Value = 269.50689
$ws.Range("O26").Value = 0.1681521229683693
$ws.Range("P26").Value = 0.1782820942415013
$ws.Range("Q26").Value = 300.43963313538
$ws.Range("R26").Value = 2703.95669821842
$ws.Range("S26").Value = 0.007414309017184866
$ws.Range("T26").Value = 0.008523525267639378
